$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new column before column C ("Verdict") - shifts old C..W right
#    (to X for the rows that already reached column W).
# ---------------------------------------------------------------------------
$ws.Columns("C:C").Insert()

# The insert drags the old column W (General numeric 0) one slot right into
# X for rows 4/5 - that spill-over column is not part of the new layout.
$ws.Columns("X:X").Delete()

# Column T (new "50Moving%") inherited the percentage number format that
# used to live one column to the left (old "S") for the rows that already
# had R:W data - put those two cells back to the default/general format.
$ws.Range("T4:T5").ClearFormats()

# ---------------------------------------------------------------------------
# 2. Header row (row 1) - full set of labels A1:W1.
# ---------------------------------------------------------------------------
$headers = @("Date","ScoreFinal","Verdict","totalSentiment","wordCount","sentenceCount","posWordPercentage","negWordPercentage","posPhrasePercentage","negPhrasePercentage","ElapsedMs","posWordCount","negWordCount","positivePhraseCount","negativePhraseCount","Method","RSI","PEG","200Moving%","50Moving%","PriceBook","Dividend","Bollinger")
for ($i = 0; $i -lt $headers.Count; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# ---------------------------------------------------------------------------
# 3. Data rows 2-5 (old rows 6 and 7 get dropped afterwards).
# ---------------------------------------------------------------------------
$dates = @(42628.850717592592, 42628.852731481478, 42628.856805555559, 42628.858796296299)
$scores = @(-12, -9, 2, 9)

for ($r = 2; $r -le 5; $r++) {
    $idx = $r - 2

    $ws.Cells.Item($r, 1).Value2 = $dates[$idx]
    $ws.Cells.Item($r, 2).Value2 = $scores[$idx]
    $ws.Cells.Item($r, 3).Value2 = "buy"

    for ($col = 4; $col -le 15; $col++) {
        $ws.Cells.Item($r, $col).Value2 = 0
    }

    $ws.Cells.Item($r, 16).Value2 = "Random"
    $ws.Cells.Item($r, 17).Value2 = 0
    $ws.Cells.Item($r, 18).Value2 = 1.66

    $c = $ws.Cells.Item($r, 19)
    $c.NumberFormat = "0.00%"
    $c.Value2 = 0.0969

    $ws.Cells.Item($r, 20).Value2 = 4.57
    $ws.Cells.Item($r, 21).Value2 = 4.5999999999999996
    $ws.Cells.Item($r, 22).Value2 = 2.2799999999999998
    $ws.Cells.Item($r, 23).Value2 = 0
}

# ---------------------------------------------------------------------------
# 4. Drop the old rows 6 and 7 (only four data rows remain now).
# ---------------------------------------------------------------------------
$ws.Rows("6:7").Delete()
